$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H112").Value = 1693.36
$ws.Range("J112").Value = 1707.75
$ws.Range("L112").Value = 5123.25
$ws.Range("N112").Value = -7339.25
$ws.Range("H121").Value = 1772.5
$ws.Range("I121").Value = 95
$ws.Range("J121").Value = 1901.5385
$ws.Range("K121").Value = 285
$ws.Range("L121").Value = 5704.6155
$ws.Range("M121").Value = 1462
$ws.Range("N121").Value = -9198.6155
$ws.Range("H129").Value = 1006.875
$ws.Range("I129").Value = 1006.875
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3020.625
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = 1979.375
$ws.Range("H131").Value = 4383.2
$ws.Range("I131").Value = 4383.2
$ws.Range("K131").Value = 13149.6
$ws.Range("M131").Value = -8109.599999999999
$ws.Range("H137").Value = 3226.6365
$ws.Range("I137").Value = 1975.5238
$ws.Range("J137").Value = 5416.0835
$ws.Range("K137").Value = 5926.5714
$ws.Range("L137").Value = 16248.2505
$ws.Range("M137").Value = -3376.5714
$ws.Range("N137").Value = -21348.2505
$ws.Range("H138").Value = 3154.175
$ws.Range("J138").Value = 5422.4443
$ws.Range("L138").Value = 16267.3329
$ws.Range("N138").Value = -26547.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H45").Value = 1054.1111
$ws.Range("I45").Value = 1142.125
$ws.Range("J45").Value = 350
$ws.Range("K45").Value = 1142.125
$ws.Range("L45").Value = 350
$ws.Range("M45").Value = -765.125
$ws.Range("N45").Value = -1104
$ws.Range("H61").Value = 2795.375
$ws.Range("I61").Value = 2766.5715
$ws.Range("K61").Value = 2766.5715
$ws.Range("M61").Value = -2554.5715
$ws.Range("H74").Value = 849.1667
$ws.Range("I74").Value = 829.55554
$ws.Range("J74").Value = 908
$ws.Range("K74").Value = 829.55554
$ws.Range("L74").Value = 908
$ws.Range("M74").Value = 44.44446000000005
$ws.Range("N74").Value = -2656
$ws.Range("H77").Value = 849.1667
$ws.Range("I77").Value = 829.55554
$ws.Range("J77").Value = 908
$ws.Range("K77").Value = 4147.7777
$ws.Range("L77").Value = 4540
$ws.Range("M77").Value = 220.2223000000004
$ws.Range("N77").Value = -13276
$ws.Range("H136").Value = 2795.375
$ws.Range("I136").Value = 2766.5715
$ws.Range("K136").Value = 8299.7145
$ws.Range("M136").Value = -5749.7145
$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H88").Value = 32500
$ws.Range("J88").Value = 32500
$ws.Range("L88").Value = 32500
$ws.Range("N88").Value = -33312
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H91").Value = 32500
$ws.Range("J91").Value = 32500
$ws.Range("L91").Value = 32500
$ws.Range("N91").Value = -35308
$ws.Range("H105").Value = 2161.2666
$ws.Range("I105").Value = 2109.2307
$ws.Range("J105").Value = 2499.5
$ws.Range("K105").Value = 2109.2307
$ws.Range("L105").Value = 2499.5
$ws.Range("M105").Value = -362.2307000000001
$ws.Range("N105").Value = -5993.5
$ws.Range("H134").Value = 2536.3333
$ws.Range("J134").Value = 2384.6667
$ws.Range("L134").Value = 7154.000100000001
$ws.Range("N134").Value = -12224.0001
$ws.Range("H135").Value = 45593.453
$ws.Range("J135").Value = 45593.453
$ws.Range("L135").Value = 45593.453
$ws.Range("N135").Value = -55733.453

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 149900
$ws.Range("J9").Value = 149900
$ws.Range("L9").Value = 149900
$ws.Range("N9").Value = -150236
$ws.Range("H122").Value = 2503.7144
$ws.Range("I122").Value = 2671
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 8013
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -5563
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3183.2856
$ws.Range("I140").Value = 3130.6667
$ws.Range("K140").Value = 9392.000100000001
$ws.Range("M140").Value = -4212.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6498
$ws.Range("J80").Value = 10006
$ws.Range("L80").Value = 10006
$ws.Range("N80").Value = -12002
$ws.Range("H83").Value = 6498
$ws.Range("J83").Value = 10006
$ws.Range("L83").Value = 50030
$ws.Range("N83").Value = -60014
$ws.Range("H122").Value = 718.6
$ws.Range("I122").Value = 717.5454999999999
$ws.Range("J122").Value = 721.5
$ws.Range("K122").Value = 2152.6365
$ws.Range("L122").Value = 2164.5
$ws.Range("M122").Value = 297.3635000000004
$ws.Range("N122").Value = -7064.5
$ws.Range("H132").Value = 3239.923
$ws.Range("I132").Value = 2902.6667
$ws.Range("J132").Value = 3998.75
$ws.Range("K132").Value = 8708.000100000001
$ws.Range("L132").Value = 11996.25
$ws.Range("M132").Value = -6178.000100000001
$ws.Range("N132").Value = -17056.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 20006
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H82").Value = 1666.6666
$ws.Range("J82").Value = 1800
$ws.Range("L82").Value = 1800
$ws.Range("N82").Value = -2522
$ws.Range("H85").Value = 1666.6666
$ws.Range("J85").Value = 1800
$ws.Range("L85").Value = 1800
$ws.Range("N85").Value = -4296

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1862.5454
$ws.Range("I132").Value = 1404.7059
$ws.Range("K132").Value = 4214.1177
$ws.Range("M132").Value = -1684.1177
$ws.Range("H136").Value = 774.3
$ws.Range("I136").Value = 809.5263
$ws.Range("J136").Value = 105
$ws.Range("K136").Value = 2428.5789
$ws.Range("L136").Value = 315
$ws.Range("M136").Value = 121.4211
$ws.Range("N136").Value = -5415
